$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "175"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "426606.00"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "426"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1389424.49"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "130"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "374445.38"

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "188"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "416926.00"

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "388"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1183532.00"

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "304"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "839218.74"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "69"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "311198.15"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "228"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "602700.74"

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "462"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1553055.70"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "319"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "987653.79"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "2197"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4792717.07"

$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "12"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "50500.00"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "3032"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8687728.66"

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "10"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44500.00"

$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "3121"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "8337519.77"

$ws.Range("C53").NumberFormat = "@"
$ws.Range("C53").Value = "43"
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = "110000.00"

$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "50"
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value = "139847.00"

$ws.Range("C66").NumberFormat = "@"
$ws.Range("C66").Value = "324"
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "777949.84"

$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value = "800"
$ws.Range("D68").NumberFormat = "@"
$ws.Range("D68").Value = "2502922.82"

$ws.Range("C69").NumberFormat = "@"
$ws.Range("C69").Value = "462"
$ws.Range("D69").NumberFormat = "@"
$ws.Range("D69").Value = "1357377.03"

$ws.Range("C71").NumberFormat = "@"
$ws.Range("C71").Value = "28"
$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = "98231.09"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "430"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "952797.75"

$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "939"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "2805543.13"

$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "857"
$ws.Range("D88").NumberFormat = "@"
$ws.Range("D88").Value = "2329163.74"

$ws.Range("C89").NumberFormat = "@"
$ws.Range("C89").Value = "12"
$ws.Range("D89").NumberFormat = "@"
$ws.Range("D89").Value = "31500.00"

